# Add files via upload
#
# The uploaded copy of this workbook has two new leading columns that were
# not present before: "col" (the collaboration/experiment, e.g. "STAR") and
# "tar" (the target/beam type, e.g. "pp"). Everything that used to live in
# columns A:O now lives in columns C:Q, and the formulas that reference
# other cells on the row follow the shift automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 7

# Insert two blank columns at the very front of the sheet; this pushes the
# existing A:O data (values, shared strings and formulas alike) to C:Q and
# keeps every formula reference consistent.
$ws.Range("A1:B1").EntireColumn.Insert()

# Header row for the two new columns.
$ws.Range("A1").Value = "col"
$ws.Range("B1").Value = "tar"

# Every data row gets the same collaboration/target values.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "STAR"
    $ws.Cells.Item($r, 2).Value = "pp"
}

# The rest of the sheet is centre-aligned (style index 1); match that for
# the newly-added columns too.
$ws.Range("A1:B" + $lastRow).HorizontalAlignment = -4108

[void]$ws.Range("A1:B" + $lastRow).Select()
